# Update "bản phân công" (task assignment) sheet:
# - Reassign testing owners in rows 5 and 7 (columns E and I)
# - Leave active cell selection on I7, matching the state the file was saved in

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5: testing owner changes from "Long" to "khang"
$ws.Range("E5").Value = "khang"
$ws.Range("I5").Value = "khang"

# Row 7: testing owner changes from "khang" to "long"
$ws.Range("E7").Value = "long"
$ws.Range("I7").Value = "long"

# Reflect the last active cell/selection at save time
[void]$ws.Range("I7").Select()
